# Week 17 data logging + Simulate_Season.py tiebreaking fix
# (the workbook itself only needs the Week-17 stat rows appended / updated)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# YDS sheet - append this week's per-play yardage logs
# ---------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$wsYDS.Range("B2").Value = ($wsYDS.Range("B2").Text + " 3 6 7 -3 5 5 4 5 2 0 1 0 6 12 4 7 4 -4 -4 2 2 0 4 5 2 5")
$wsYDS.Range("B3").Value = ($wsYDS.Range("B3").Text + " 15 0 44 18 -1 24 8 35 7 9 5 11 6 2 16 7 34 5")
$wsYDS.Range("C2").Value = ($wsYDS.Range("C2").Text + " 4 9 -1 2 14 3 2 6 4 3 1 6 9 6 0 2 0 -1 1 2 -4 2 7 4 5 -4 2 5 7 2 4 1 3")
$wsYDS.Range("C3").Value = ($wsYDS.Range("C3").Text + " 10 8 10 11 8 3 6 8 3 5 9 3 8 7 40 16 5 5 45 8 8")

# ---------------------------------------------------------------
# OFF sheet - season totals through Week 17 (row2 = RATT, row3 = PATT)
# ---------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Range("B2").Value = 7
$wsOFF.Range("C2").Value = 194
$wsOFF.Range("E2").Value = 18
$wsOFF.Range("F2").Value = 64
$wsOFF.Range("G2").Value = 59
$wsOFF.Range("I2").Value = 9
$wsOFF.Range("J2").Value = 30
$wsOFF.Range("N2").Value = 17
$wsOFF.Range("O2").Value = 28
$wsOFF.Range("P2").Value = 16

$wsOFF.Range("C3").Value = 165
$wsOFF.Range("D3").Value = 6
$wsOFF.Range("E3").Value = 34
$wsOFF.Range("F3").Value = 92
$wsOFF.Range("G3").Value = 28
$wsOFF.Range("I3").Value = 70
$wsOFF.Range("J3").Value = 40
$wsOFF.Range("L3").Value = 262
$wsOFF.Range("M3").Value = 176
$wsOFF.Range("Q3").Value = 488

# ---------------------------------------------------------------
# DEF sheet - season totals through Week 17 (row2 = RATT, row3 = PATT)
# ---------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Range("C2").Value = 191
$wsDEF.Range("D2").Value = 7
$wsDEF.Range("F2").Value = 68
$wsDEF.Range("G2").Value = 52
$wsDEF.Range("J2").Value = 38
$wsDEF.Range("O2").Value = 25
$wsDEF.Range("P2").Value = 10

$wsDEF.Range("C3").Value = 148
$wsDEF.Range("E3").Value = 31
$wsDEF.Range("F3").Value = 85
$wsDEF.Range("G3").Value = 43
$wsDEF.Range("I3").Value = 51
$wsDEF.Range("J3").Value = 64
$wsDEF.Range("L3").Value = 251
$wsDEF.Range("M3").Value = 148
$wsDEF.Range("Q3").Value = 476

# ---------------------------------------------------------------
# ST sheet - kicking / punting season totals + per-kick logs
# ---------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Range("B2").Value = 72
$wsST.Range("D2").Value = 62
$wsST.Range("F2").Value = 256
$wsST.Range("G2").Value = 249
$wsST.Range("J2").Value = 117
$wsST.Range("K2").Value = 111
$wsST.Range("P2").Value = 2
$wsST.Range("Q2").Value = 1

$wsST.Range("B3").Value = 58
$wsST.Range("D3").Value = ($wsST.Range("D3").Text + " 41 56 40")

$wsST.Range("B4").Value = ($wsST.Range("B4").Text + " 67 66")
$wsST.Range("D4").Value = ($wsST.Range("D4").Text + " 0 1 4")

$wsST.Range("B5").Value = ($wsST.Range("B5").Text + " 47 1")
$wsST.Range("D5").Value = ($wsST.Range("D5").Text + " 0 4")

# ---------------------------------------------------------------
# TURNS sheet - fumbles lost corrected
# ---------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Range("E3").Value = 5

# ---------------------------------------------------------------
# PEN sheet - penalty counts
# ---------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Range("B2").Value = 16
$wsPEN.Range("D4").Value = 5
